$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2.86
$ws.Range("H3").Value = 2.42
$ws.Range("I3").Value = 2.64
$ws.Range("J3").Value = 3.25
$ws.Range("P3").Value = 1.95
$ws.Range("Q3").Value = 1.89
$ws.Range("H4").Value = 2.52
$ws.Range("I4").Value = 2.82
$ws.Range("K4").Value = 3.2
$ws.Range("N4").Value = 2.58
$ws.Range("O4").Value = 1.52
$ws.Range("X4").Value = 9.199999999999999
$ws.Range("AA4").Value = 48
$ws.Range("AE4").Value = 44
$ws.Range("AI4").Value = 70
$ws.Range("AK4").Value = 60
$ws.Range("AM4").Value = 190
$ws.Range("AO4").Value = 46
$ws.Range("F6").Value = 3.2
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 2.28
$ws.Range("I6").Value = 2.56
$ws.Range("J6").Value = 3.15
$ws.Range("K6").Value = 3.7
$ws.Range("F8").Value = 5.8
$ws.Range("I9").Value = 2.24
$ws.Range("F10").Value = 3.45
$ws.Range("G10").Value = 3.6
$ws.Range("J10").Value = 3.55
$ws.Range("K10").Value = 3.65
$ws.Range("AF10").Value = 26
$ws.Range("F11").Value = 4.3
$ws.Range("I11").Value = 2
$ws.Range("K11").Value = 3.75
$ws.Range("P11").Value = 1.9
$ws.Range("X11").Value = 1000
$ws.Range("Y11").Value = 10
$ws.Range("Z11").Value = 1000
$ws.Range("AA11").Value = 980
$ws.Range("AB11").Value = 1000
$ws.Range("AC11").Value = 1000
$ws.Range("AD11").Value = 1000
$ws.Range("AE11").Value = 980
$ws.Range("AF11").Value = 980
$ws.Range("AG11").Value = 1000
$ws.Range("AH11").Value = 1000
$ws.Range("AI11").Value = 980
$ws.Range("AJ11").Value = 110
$ws.Range("AO11").Value = 1000
$ws.Range("F12").Value = 1.33
$ws.Range("G12").Value = 1.35
$ws.Range("H12").Value = 11
$ws.Range("I12").Value = 13
$ws.Range("J12").Value = 5.7
$ws.Range("K12").Value = 6.4
$ws.Range("Q12").Value = 1.61
$ws.Range("Y12").Value = 44
$ws.Range("Z12").Value = 130
$ws.Range("AA12").Value = 540
$ws.Range("AD12").Value = 50
$ws.Range("AH12").Value = 34
$ws.Range("AL12").Value = 44
$ws.Range("Q13").Value = 1.84
$ws.Range("T13").Value = 1.85
$ws.Range("X13").Value = 17.5
$ws.Range("AD13").Value = 24
$ws.Range("AJ13").Value = 20
$ws.Range("G15").Value = 7
$ws.Range("K15").Value = 5.1
$ws.Range("N15").Value = 5
$ws.Range("P15").Value = 2.42
$ws.Range("R15").Value = 1.57
$ws.Range("S15").Value = 2.5
$ws.Range("T15").Value = 1.76
$ws.Range("U15").Value = 2.16
$ws.Range("AJ15").Value = 200
$ws.Range("P16").Value = 1.78
$ws.Range("F17").Value = 1.46
$ws.Range("H17").Value = 8.4
$ws.Range("I17").Value = 9.199999999999999
$ws.Range("H19").Value = 2.62
$ws.Range("J19").Value = 3.25
$ws.Range("P20").Value = 1.72
